$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows after the current row 219 (existing 220 -> 222,
# two fresh rows become 220 and 221).
$ws.Rows("220:221").Insert()

# --- Row 219: update in place with the new weekly reading ---
$ws.Cells.Item(219, 4).Value = 44595              # D: Fecha
$ws.Cells.Item(219, 10).Value = 90                # J: Volumen
$ws.Cells.Item(219, 11).Value = 7000              # K: Precio minimo (unchanged)
$ws.Cells.Item(219, 12).Value = 7000              # L: Precio maximo
$ws.Cells.Item(219, 13).Value = 7000              # M: Precio promedio ponderado
$ws.Cells.Item(219, 15).Value = "Provincia del Elquí"   # O: Origen
$ws.Cells.Item(219, 16).Value = 1167              # P: Precio $/Kg

# --- Row 220 (newly inserted): another new weekly reading ---
$ws.Cells.Item(220, 1).Value = 9
$ws.Cells.Item(220, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(220, 3).Value = "Metropolitana"
$ws.Cells.Item(220, 4).Value = 44595
$ws.Cells.Item(220, 5).Value = 13
$ws.Cells.Item(220, 6).Value = 100112017
$ws.Cells.Item(220, 7).Value = "Apio"
$ws.Cells.Item(220, 8).Value = "Americana (o)"
$ws.Cells.Item(220, 9).Value = "Segunda"
$ws.Cells.Item(220, 10).Value = 50
$ws.Cells.Item(220, 11).Value = 5000
$ws.Cells.Item(220, 12).Value = 5000
$ws.Cells.Item(220, 13).Value = 5000
$ws.Cells.Item(220, 14).Value = "`$/docena de matas"
$ws.Cells.Item(220, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(220, 16).Value = 833
$ws.Cells.Item(220, 17).Value = 6
$ws.Cells.Item(220, 18).Value = "Hortaliza"

# --- Row 221 (newly inserted): carries the previous (pre-edit) row 219 reading ---
$ws.Cells.Item(221, 1).Value = 9
$ws.Cells.Item(221, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(221, 3).Value = "Metropolitana"
$ws.Cells.Item(221, 4).Value = 44544
$ws.Cells.Item(221, 5).Value = 13
$ws.Cells.Item(221, 6).Value = 100112017
$ws.Cells.Item(221, 7).Value = "Apio"
$ws.Cells.Item(221, 8).Value = "Americana (o)"
$ws.Cells.Item(221, 9).Value = "Primera"
$ws.Cells.Item(221, 10).Value = 79
$ws.Cells.Item(221, 11).Value = 7000
$ws.Cells.Item(221, 12).Value = 8000
$ws.Cells.Item(221, 13).Value = 7494
$ws.Cells.Item(221, 14).Value = "`$/docena de matas"
$ws.Cells.Item(221, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(221, 16).Value = 1249
$ws.Cells.Item(221, 17).Value = 6
$ws.Cells.Item(221, 18).Value = "Hortaliza"

# Row 222 keeps the original (pre-edit) row 220 content automatically, since
# the Insert() above shifted it down without altering its values.

# Make sure the D-column date cells keep the date number format (style index 2)
$dateFmt = $ws.Cells.Item(218, 4).NumberFormat
$ws.Cells.Item(219, 4).NumberFormat = $dateFmt
$ws.Cells.Item(220, 4).NumberFormat = $dateFmt
$ws.Cells.Item(221, 4).NumberFormat = $dateFmt
$ws.Cells.Item(222, 4).NumberFormat = $dateFmt
